# Updated cryptos list on Sun Sep 10 10:53:57 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# swaps the #10/#11 ranked coins (WrappedEther <-> WrappedliquidstakedEther2.0)
# including their link + price + volume cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that render as plain numbers (e.g. "213.80", "0.0623")
# get silently coerced to a Double by COM's Range.Value setter, which loses
# the exact text (trailing zeros, leading zero digits, decimal grouping).
# Forcing a "@" (text) number format before the write keeps the literal
# string; ClearFormats() afterwards drops the now-unneeded text format so
# the cell's style stays identical to every other untouched data cell.

$ws.Range('D2').Value = '25.985.20'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.632.87'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.80'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.252'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0623'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.49'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -5.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0790'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '1.860.54'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.660.53'
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.19'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('D16').Value = '25.994.64'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('E17').Value = '  -2.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.73'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.88%  '
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '190.32'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.24'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.56'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.12'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.07%  '
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.61'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.76'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.21'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.09%  '
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0483'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.11%  '
$ws.Range('E32').Value = '  -3.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.14'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.28%  '
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.49'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.34%  '
$ws.Range('D36').Value = '1.134.02'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.866'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.14%  '
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('E39').Value = '  -2.98%  '
$ws.Range('E40').Value = '  -1.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.56'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.780'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('E43').Value = '  -4.78%  '
$ws.Range('D44').Value = '1.770.49'
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.08'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.73%  '
$ws.Range('E47').Value = '  -0.70%  '
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.51'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.98%  '
$ws.Range('E51').Value = '  +0.33%  '
